# Posicion_Balizas.xlsx - "Added integration project and small modifications
# to LM positions"
#
# The LM 1..4 X/Y positions (B3:C6) change value. In the source workbook
# these numeric-looking values are stored as shared-string TEXT (not
# numbers) using the same cell style as every other data cell (s="1", no
# quote-prefix). Writing the value directly (Range.Value = "0.1") would be
# auto-coerced to a real number by Excel, and forcing text via a leading
# apostrophe (or NumberFormat "@") stamps a new quotePrefix style onto the
# cell. Instead we compute the text through TEXT() and then overwrite the
# formula with its own computed value (paste-values-only), which leaves a
# plain text cell behind using the cell's existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber($cellRef, $numericText) {
    $cell = $ws.Range($cellRef)
    $cell.Formula = '=TEXT(' + $numericText + ',"0.0")'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# LM 1: (row 3)
Set-TextNumber "B3" "0.1"
Set-TextNumber "C3" "0.1"

# LM 2: (row 4)
Set-TextNumber "B4" "0.1"
Set-TextNumber "C4" "13.4"

# LM 3: (row 5)
Set-TextNumber "B5" "26.4"
Set-TextNumber "C5" "13.4"

# LM 4: (row 6)
Set-TextNumber "B6" "26.4"
Set-TextNumber "C6" "0.1"

$excel.CutCopyMode = $false

# Cursor/selection left on I8, matching the saved sheet view.
$ws.Range("I8").Select() | Out-Null
